# Updates cryptos list prices/volume figures on the active worksheet.
# Values in column D that look like plain decimal numbers are written with a
# leading apostrophe so Excel keeps storing them as text (matching the
# original workbook, where every Price/Volume cell is text, not a number).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '61.332.22'
$ws.Range('E2').Value = '  +0.64%  '
$ws.Range('D3').Value = '2.923.82'
$ws.Range('E3').Value = '  +0.11%  '
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').Value = '''597.68'
$ws.Range('E5').Value = '  +0.66%  '
$ws.Range('D6').Value = '''144.91'
$ws.Range('E6').Value = '  -0.53%  '
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('E8').Value = '  -1.17%  '
$ws.Range('D9').Value = '''6.92'
$ws.Range('E9').Value = '  +1.19%  '
$ws.Range('E10').Value = '  -1.90%  '
$ws.Range('E11').Value = '  -0.48%  '
$ws.Range('E12').Value = '  -0.90%  '
$ws.Range('D13').Value = '''33.45'
$ws.Range('E13').Value = '  -0.68%  '
$ws.Range('E14').Value = '  +0.01%  '
$ws.Range('D15').Value = '3.409.52'
$ws.Range('E15').Value = '  +0.21%  '
$ws.Range('D16').Value = '61.241.20'
$ws.Range('E16').Value = '  +0.53%  '
$ws.Range('D17').Value = '2.926.00'
$ws.Range('E17').Value = '  +0.14%  '
$ws.Range('E18').Value = '  -0.54%  '
$ws.Range('D19').Value = '''431.84'
$ws.Range('E19').Value = '  +0.23%  '
$ws.Range('E20').Value = '  +1.38%  '
$ws.Range('D21').Value = '''0.674'
$ws.Range('E21').Value = '  -1.24%  '
$ws.Range('E22').Value = '  -0.04%  '
$ws.Range('D23').Value = '''81.82'
$ws.Range('E23').Value = '  +0.41%  '
$ws.Range('E24').Value = '  -0.56%  '
$ws.Range('E25').Value = '  -1.43%  '
$ws.Range('D26').Value = '''11.75'
$ws.Range('E26').Value = '  -1.68%  '
$ws.Range('E27').Value = '  +0.00%  '
$ws.Range('E28').Value = '  -4.42%  '
$ws.Range('E29').Value = '  -0.61%  '
$ws.Range('D30').Value = '''6.89'
$ws.Range('E30').Value = '  -2.48%  '
$ws.Range('D31').Value = '''26.62'
$ws.Range('E31').Value = '  +0.62%  '
$ws.Range('E32').Value = '  +1.40%  '
$ws.Range('E33').Value = '  +0.08%  '
$ws.Range('D34').Value = '0.0₃0878'
$ws.Range('E34').Value = '  +3.40%  '
$ws.Range('E35').Value = '  -0.32%  '
$ws.Range('E36').Value = '  -0.10%  '
$ws.Range('E37').Value = '  -1.29%  '
$ws.Range('E38').Value = '  +0.60%  '
$ws.Range('E39').Value = '  -0.88%  '
$ws.Range('E40').Value = '  -0.28%  '
$ws.Range('D41').Value = '''42.48'
$ws.Range('E41').Value = '  +4.67%  '
$ws.Range('D42').Value = '''0.280'
$ws.Range('E42').Value = '  -2.22%  '
$ws.Range('E43').Value = '  -0.41%  '
$ws.Range('D44').Value = '2.694.82'
$ws.Range('E44').Value = '  -0.66%  '
$ws.Range('D45').Value = '''366.23'
$ws.Range('E45').Value = '  -2.14%  '
$ws.Range('D46').Value = '''133.59'
$ws.Range('E46').Value = '  +2.25%  '
$ws.Range('D48').Value = '''23.50'
$ws.Range('E48').Value = '  -2.10%  '
$ws.Range('E49').Value = '  -1.25%  '
$ws.Range('E50').Value = '  -0.99%  '
$ws.Range('E51').Value = '  -1.00%  '
